# Auto update Excel log
# Appends new sensor/alert rows to the ALERTS, PIR, and Humidity sheets.
# NumberFormat is forced to Text ("@") on the new rows before assigning
# values so date-looking ("2026-01-30") and percent-looking ("87.1%")
# strings are preserved verbatim instead of being coerced by Excel into
# date serials / percentage numbers.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# ALERTS sheet: two new FALL_DETECTED rows (15-16)
# ---------------------------------------------------------------------------
$wsAlerts = $wb.Worksheets.Item("ALERTS")
$alertsRows = @(
    ,@("2026-01-30", "17:50:45", "17:00", "Living Room", "CRITICAL", "FALL_DETECTED")
    ,@("2026-01-30", "17:50:48", "17:00", "Living Room", "CRITICAL", "FALL_DETECTED")
)
$wsAlerts.Range("A15:F16").NumberFormat = "@"
$r = 15
foreach ($row in $alertsRows) {
    $wsAlerts.Cells.Item($r, 1).Value = $row[0]
    $wsAlerts.Cells.Item($r, 2).Value = $row[1]
    $wsAlerts.Cells.Item($r, 3).Value = $row[2]
    $wsAlerts.Cells.Item($r, 4).Value = $row[3]
    $wsAlerts.Cells.Item($r, 5).Value = $row[4]
    $wsAlerts.Cells.Item($r, 6).Value = $row[5]
    $r++
}

# ---------------------------------------------------------------------------
# PIR sheet: eleven new "No Motion" / Inactive rows (346-356)
# ---------------------------------------------------------------------------
$wsPir = $wb.Worksheets.Item("PIR")
$pirRows = @(
    ,@("2026-01-30", "17:50:51", "17:00", "Bathroom", "No Motion", "Inactive")
    ,@("2026-01-30", "17:50:56", "17:00", "Bathroom", "No Motion", "Inactive")
    ,@("2026-01-30", "17:51:01", "17:00", "Bathroom", "No Motion", "Inactive")
    ,@("2026-01-30", "17:51:06", "17:00", "Bathroom", "No Motion", "Inactive")
    ,@("2026-01-30", "17:51:11", "17:00", "Bathroom", "No Motion", "Inactive")
    ,@("2026-01-30", "17:51:16", "17:00", "Bathroom", "No Motion", "Inactive")
    ,@("2026-01-30", "17:51:21", "17:00", "Bathroom", "No Motion", "Inactive")
    ,@("2026-01-30", "17:51:26", "17:00", "Bathroom", "No Motion", "Inactive")
    ,@("2026-01-30", "17:51:31", "17:00", "Bathroom", "No Motion", "Inactive")
    ,@("2026-01-30", "17:51:36", "17:00", "Bathroom", "No Motion", "Inactive")
    ,@("2026-01-30", "17:51:41", "17:00", "Bathroom", "No Motion", "Inactive")
)
$wsPir.Range("A346:F356").NumberFormat = "@"
$r = 346
foreach ($row in $pirRows) {
    $wsPir.Cells.Item($r, 1).Value = $row[0]
    $wsPir.Cells.Item($r, 2).Value = $row[1]
    $wsPir.Cells.Item($r, 3).Value = $row[2]
    $wsPir.Cells.Item($r, 4).Value = $row[3]
    $wsPir.Cells.Item($r, 5).Value = $row[4]
    $wsPir.Cells.Item($r, 6).Value = $row[5]
    $r++
}

# ---------------------------------------------------------------------------
# Humidity sheet: five new "Active" rows (239-243)
# ---------------------------------------------------------------------------
$wsHumidity = $wb.Worksheets.Item("Humidity")
$humidityRows = @(
    ,@("2026-01-30", "17:50:57", "17:00", "Bathroom", "87.1%", "Active")
    ,@("2026-01-30", "17:51:07", "17:00", "Bathroom", "86.1%", "Active")
    ,@("2026-01-30", "17:51:17", "17:00", "Bathroom", "87.1%", "Active")
    ,@("2026-01-30", "17:51:22", "17:00", "Bathroom", "87.0%", "Active")
    ,@("2026-01-30", "17:51:38", "17:00", "Bathroom", "87.0%", "Active")
)
$wsHumidity.Range("A239:F243").NumberFormat = "@"
$r = 239
foreach ($row in $humidityRows) {
    $wsHumidity.Cells.Item($r, 1).Value = $row[0]
    $wsHumidity.Cells.Item($r, 2).Value = $row[1]
    $wsHumidity.Cells.Item($r, 3).Value = $row[2]
    $wsHumidity.Cells.Item($r, 4).Value = $row[3]
    $wsHumidity.Cells.Item($r, 5).Value = $row[4]
    $wsHumidity.Cells.Item($r, 6).Value = $row[5]
    $r++
}
